$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.307.84"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "1.876.45"
$ws.Range("E3").Value = "  +1.05%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("D6").Value = "242.32"
$ws.Range("E6").Value = "  +0.84%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "0.3109"
$ws.Range("E8").Value = "  +1.24%  "

$ws.Range("D9").Value = "0.07757"
$ws.Range("E9").Value = "  +0.19%  "

$ws.Range("D10").Value = "24.96"
$ws.Range("E10").Value = "  -0.71%  "

$ws.Range("D11").Value = "0.08536"
$ws.Range("E11").Value = "  +3.48%  "

$ws.Range("D12").Value = "1.891.04"
$ws.Range("E12").Value = "  +3.06%  "

$ws.Range("E13").Value = "  -0.44%  "

$ws.Range("D14").Value = "0.7106"
$ws.Range("E14").Value = "  -0.78%  "

$ws.Range("D15").Value = "91.44"
$ws.Range("E15").Value = "  +1.38%  "

$ws.Range("D16").Value = "29.303.70"
$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("E17").Value = "  +5.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.010"
$ws.Range("E18").Value = "  +2.48%  "

$ws.Range("D19").Value = "242.27"
$ws.Range("E19").Value = "  -0.56%  "

$ws.Range("D20").Value = "2.135.23"
$ws.Range("E20").Value = "  +1.37%  "

$ws.Range("D21").Value = "13.25"
$ws.Range("E21").Value = "  +0.77%  "

$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").Value = "7.824"
$ws.Range("E23").Value = "  -2.14%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").Value = "0.1608"
$ws.Range("E25").Value = "  +0.84%  "

$ws.Range("D26").Value = "162.92"
$ws.Range("E26").Value = "  +0.30%  "

$ws.Range("D27").Value = "9.044"
$ws.Range("E27").Value = "  +1.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.50"
$ws.Range("E28").Value = "  +1.25%  "

$ws.Range("E29").Value = "  +1.13%  "

$ws.Range("D30").Value = "4.403"
$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("D31").Value = "4.317"
$ws.Range("E31").Value = "  +2.98%  "

$ws.Range("E32").Value = "  -2.70%  "

$ws.Range("D33").Value = "0.05249"
$ws.Range("E33").Value = "  +1.23%  "

$ws.Range("D35").Value = "1.176"
$ws.Range("E35").Value = "  +0.41%  "

$ws.Range("D36").Value = "0.7466"
$ws.Range("E36").Value = "  +2.74%  "

$ws.Range("D37").Value = "2.686"
$ws.Range("E37").Value = "  +0.33%  "

$ws.Range("D38").Value = "0.01867"
$ws.Range("E38").Value = "  +0.70%  "

$ws.Range("D39").Value = "2.718"
$ws.Range("E39").Value = "  +1.34%  "

$ws.Range("D40").Value = "1.183.45"
$ws.Range("E40").Value = "  +2.34%  "

$ws.Range("D41").Value = "6.392"
$ws.Range("E41").Value = "  +4.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.50"
$ws.Range("E44").Value = "  +4.89%  "

$ws.Range("D46").Value = "2.030.72"
$ws.Range("E46").Value = "  +1.31%  "

$ws.Range("E47").Value = "  +2.61%  "

$ws.Range("E48").Value = "  -0.17%  "

$ws.Range("E49").Value = "  +1.55%  "

$ws.Range("D50").Value = "9.403"
$ws.Range("E50").Value = "  +0.95%  "

$ws.Range("D51").Value = "0.4315"
$ws.Range("E51").Value = "  +1.28%  "

# Row 42 becomes Aave, Row 43 becomes TrustWalletToken (content swap)
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "73.03"
$ws.Range("E42").Value = "  +1.18%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.8878"
$ws.Range("E43").Value = "  -1.53%  "
